# Add a new item row ("كريم ONE للبشره الحساسه") as item #16, just above the
# totals row, and bump the total accordingly (timestamp bump is implicit in
# the new save).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 19 - this pushes the existing totals
# row (19 -> 20) and the footer row (20 -> 21) down by one, shifting any
# merged ranges below it automatically.
$ws.Rows.Item(19).Insert()

# Clone the formatting of the row above (row 18, the last product row) onto
# the freshly inserted row 19 so it looks like the other data rows.
$ws.Range("A18:N18").Copy()
$ws.Range("A19:N19").PasteSpecial(-4122)

# Recreate the same merge pattern used by every other data row.
$ws.Range("B19:G19").Merge()
$ws.Range("H19:K19").Merge()
$ws.Range("L19:M19").Merge()

# Populate the new row with the new product's data (item #16).
$ws.Cells.Item(19, 1).Value = 16
$ws.Cells.Item(19, 2).Value = "كريم ONE للبشره الحساسه"
$ws.Cells.Item(19, 8).Value = "20:0"
$ws.Cells.Item(19, 12).Value = 25
$ws.Cells.Item(19, 14).Value = "1:0"

# Match the row heights from the target layout.
$ws.Rows.Item(19).RowHeight = 24.75
$ws.Rows.Item(20).RowHeight = 26.25

# Update the grand total (old row 19, now row 20) to include the new row.
$ws.Cells.Item(20, 11).Value = 578.32
